$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new row 19 data
$ws.Range("A19").Value = "F3"
$ws.Range("B19").Value = 20
$ws.Range("C19").Value = "YİYECEK"

# Match style (centered alignment) used by the rest of the table
$ws.Range("A19:C19").HorizontalAlignment = -4108

# Update the active selection to match the new cell position
$ws.Range("C19").Select()
